# Updated project: Glaucoma test app with full features
# Applies the new Pelli-Robinson contrast sensitivity test entry (row 4)
# and removes the stray empty Doctor Notes cell (I3) on that sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pelli-Robinson")

# Remove the empty inline-string placeholder cell at I3 so it no longer
# exists in the sheet (was previously an empty inlineStr cell).
$ws.Range("I3").ClearContents()

# Add the new test-result row (row 4).
$ws.Range("A4").Value = "ANI"

# Force B4 to be stored as literal text (not auto-converted to a date
# serial number), then reset its style so no extra formatting sticks.
$ws.Range("B4").NumberFormat = "@"
$ws.Range("B4").Value = "2025-05-27"
$ws.Range("B4").Style = "Normal"

$ws.Range("C4").Value = "8:51:46 PM"
$ws.Range("D4").Value = "8:52:40 PM"
$ws.Range("E4").Value = 54
$ws.Range("F4").Value = 72
$ws.Range("G4").Value = 72
$ws.Range("H4").Value = 100
$ws.Range("J4").Value = "english"
$ws.Range("K4").Value = 1.15
$ws.Range("L4").Value = 1.15
